# EditProfileName Test Case Added
#
# Inserts three new leading columns (FirstName, LastName, FullName) on the
# "Profile" worksheet, reusing the existing SignUp name/email data already
# present elsewhere in the workbook ("Sheila" / "Dimasuhid") and adding the
# concatenated "Sheila Dimasuhid" full name, then makes "Profile" the active
# sheet/tab (previously "ShareSkill" was active).

$wb = $excel.ActiveWorkbook

$wsProfile = $wb.Worksheets.Item("Profile")

# Shift the existing Profile columns (A:J) three columns to the right so the
# new name columns can be inserted ahead of them (A:C), matching the diff's
# D:M destination range for the original headers/data.
$wsProfile.Columns("A:C").Insert()

# New header row (row 1).
$wsProfile.Range("A1").Value = "FirstName"
$wsProfile.Range("B1").Value = "LastName"
$wsProfile.Range("C1").Value = "FullName"

# New data row (row 2).
$wsProfile.Range("A2").Value = "Sheila"
$wsProfile.Range("B2").Value = "Dimasuhid"
$wsProfile.Range("C2").Value = "Sheila Dimasuhid"

# Give the new columns explicit widths approximating the recorded best-fit
# widths for their header/value content.
$wsProfile.Columns("A").ColumnWidth = 8.5
$wsProfile.Columns("B").ColumnWidth = 18.5
$wsProfile.Columns("C").ColumnWidth = 14.3

# Make "Profile" the active sheet/tab with B3 selected. This also clears
# tabSelected on the previously-active "ShareSkill" sheet.
[void]$wsProfile.Activate()
[void]$wsProfile.Range("B3").Select()
